# Word COM-interop script reproducing the target edit:
#  1. Split the single run "Inside file" into three runs:
#     "Inside " / "DOCX " / "file" (and drop the stray empty <w:rPr/>
#     that sat directly under the paragraph's <w:pPr/>).
#  2. Update the "Normal" paragraph style:
#       - font color: auto -> 00000A
#       - paragraph reading order (bidi): explicit left-to-right (0)
#       - paragraph alignment: explicit left (wdAlignParagraphLeft)

$d = $word.ActiveDocument

# --- 1. Rebuild the first paragraph's runs -------------------------------
$para = $d.Paragraphs(1).Range

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:pPr><w:pStyle w:val="Normal"/></w:pPr>' +
              '<w:r><w:rPr/><w:t xml:space="preserve">Inside </w:t></w:r>' +
              '<w:r><w:rPr/><w:t xml:space="preserve">DOCX </w:t></w:r>' +
              '<w:r><w:rPr/><w:t>file</w:t></w:r>' +
              '<w:r/>' +
              '</w:p>'

$para.InsertXML($newParaXml)

# --- 2. Tweak the "Normal" style ------------------------------------------
$normal = $d.Styles("Normal")

# w:color w:val="00000A" (Word color longs are 0x00BBGGRR)
$normal.Font.Color = 0xA0000

# <w:bidi w:val="0"/>
$normal.ParagraphFormat.ReadingOrder = 0

# <w:jc w:val="left"/>
$normal.ParagraphFormat.Alignment = 0
